$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.35
$ws.Range("C2").Value = 2.2
$ws.Range("D2").Value = 7.4
$ws.Range("E2").Value = 42.25
$ws.Range("F2").Value = 26.75

$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 4.4
$ws.Range("D3").Value = 15.2
$ws.Range("E3").Value = 87.84999999999999
$ws.Range("F3").Value = 57.05

$ws.Range("B4").Value = 10.2
$ws.Range("C4").Value = 7.45
$ws.Range("D4").Value = 22.15
$ws.Range("E4").Value = 132.9
$ws.Range("F4").Value = 87.65000000000001

$ws.Range("B5").Value = 14.2
$ws.Range("C5").Value = 9.550000000000001
$ws.Range("D5").Value = 31.35
$ws.Range("E5").Value = 176.1
$ws.Range("F5").Value = 119.05

$ws.Range("B6").Value = 16.65
$ws.Range("C6").Value = 13.1
$ws.Range("D6").Value = 38.55
$ws.Range("E6").Value = 216.3
$ws.Range("F6").Value = 152.05

$ws.Range("B7").Value = 21.5
$ws.Range("C7").Value = 14.8
$ws.Range("D7").Value = 48.4
$ws.Range("E7").Value = 273.25
$ws.Range("F7").Value = 186.15

$ws.Range("B8").Value = 25.2
$ws.Range("C8").Value = 17.5
$ws.Range("D8").Value = 59.25
$ws.Range("E8").Value = 318.1
$ws.Range("F8").Value = 218.7

$ws.Range("B9").Value = 29.05
$ws.Range("C9").Value = 19.75
$ws.Range("D9").Value = 66.2
$ws.Range("E9").Value = 360.55
$ws.Range("F9").Value = 253.6

$ws.Range("B10").Value = 32.6
$ws.Range("C10").Value = 22.35
$ws.Range("D10").Value = 75.09999999999999
$ws.Range("E10").Value = 402.85
$ws.Range("F10").Value = 289.1

$ws.Range("B11").Value = 35.15
$ws.Range("C11").Value = 27.35
$ws.Range("D11").Value = 81.7
$ws.Range("E11").Value = 441.8
$ws.Range("F11").Value = 320.75

$ws.Range("B12").Value = 40.2
$ws.Range("C12").Value = 29.9
$ws.Range("D12").Value = 88.55
$ws.Range("E12").Value = 499.25
$ws.Range("F12").Value = 356.05

$ws.Range("B13").Value = 42.45
$ws.Range("C13").Value = 33
$ws.Range("D13").Value = 99.09999999999999
$ws.Range("E13").Value = 558.7
$ws.Range("F13").Value = 390.75

$ws.Range("B14").Value = 47.7
$ws.Range("C14").Value = 33.9
$ws.Range("D14").Value = 115.05
$ws.Range("E14").Value = 604.75
$ws.Range("F14").Value = 424.8

$ws.Range("B15").Value = 52.1
$ws.Range("C15").Value = 37
$ws.Range("D15").Value = 123.9
$ws.Range("E15").Value = 647.9
$ws.Range("F15").Value = 460.45

$ws.Range("B16").Value = 60
$ws.Range("C16").Value = 40.6
$ws.Range("D16").Value = 135.05
$ws.Range("E16").Value = 708.85
$ws.Range("F16").Value = 502.25

$ws.Range("B17").Value = 63.65
$ws.Range("C17").Value = 43.1
$ws.Range("D17").Value = 143.5
$ws.Range("E17").Value = 753.5
$ws.Range("F17").Value = 536.7

$ws.Range("B18").Value = 65
$ws.Range("C18").Value = 44.4
$ws.Range("D18").Value = 148.65
$ws.Range("E18").Value = 784
$ws.Range("F18").Value = 570.15

$ws.Range("B19").Value = 66.75
$ws.Range("C19").Value = 45.85
$ws.Range("D19").Value = 155.65
$ws.Range("E19").Value = 820.45
$ws.Range("F19").Value = 599.7

$ws.Range("B20").Value = 69.7
$ws.Range("C20").Value = 51.95
$ws.Range("D20").Value = 164.55
$ws.Range("E20").Value = 864.25
$ws.Range("F20").Value = 638.75

$ws.Range("B21").Value = 73.09999999999999
$ws.Range("C21").Value = 59.3
$ws.Range("D21").Value = 173.2
$ws.Range("E21").Value = 904.65
$ws.Range("F21").Value = 676.5

$ws.Range("B22").Value = 74.3
$ws.Range("C22").Value = 59.5
$ws.Range("D22").Value = 183.8
$ws.Range("E22").Value = 958.4
$ws.Range("F22").Value = 712.95

$ws.Range("B23").Value = 81.05
$ws.Range("C23").Value = 63.6
$ws.Range("D23").Value = 189.4
$ws.Range("E23").Value = 1031
$ws.Range("F23").Value = 750.65

$ws.Range("B24").Value = 85.15000000000001
$ws.Range("C24").Value = 65.34999999999999
$ws.Range("D24").Value = 196.9
$ws.Range("E24").Value = 1092.7
$ws.Range("F24").Value = 790.35

$ws.Range("B25").Value = 91.59999999999999
$ws.Range("C25").Value = 69
$ws.Range("D25").Value = 212.75
$ws.Range("E25").Value = 1157.85
$ws.Range("F25").Value = 829.55

$ws.Range("B26").Value = 95.59999999999999
$ws.Range("C26").Value = 70.5
$ws.Range("D26").Value = 235.2
$ws.Range("E26").Value = 1194.4
$ws.Range("F26").Value = 857.6

$ws.Range("B27").Value = 99
$ws.Range("C27").Value = 71.95
$ws.Range("D27").Value = 244.1
$ws.Range("E27").Value = 1244.15
$ws.Range("F27").Value = 897.2

$ws.Range("B28").Value = 105.6
$ws.Range("C28").Value = 75.65000000000001
$ws.Range("D28").Value = 253.35
$ws.Range("E28").Value = 1292.05
$ws.Range("F28").Value = 931.5

$ws.Range("B29").Value = 111.75
$ws.Range("C29").Value = 78.7
$ws.Range("D29").Value = 261.8
$ws.Range("E29").Value = 1343.1
$ws.Range("F29").Value = 974.6

$ws.Range("B30").Value = 115.5
$ws.Range("C30").Value = 82.2
$ws.Range("D30").Value = 274.5
$ws.Range("E30").Value = 1396.1
$ws.Range("F30").Value = 1016.5

$ws.Range("B31").Value = 120.7
$ws.Range("C31").Value = 82.65000000000001
$ws.Range("D31").Value = 284.1
$ws.Range("E31").Value = 1447.8
$ws.Range("F31").Value = 1052.85

$ws.Range("B32").Value = 126.4
$ws.Range("C32").Value = 86.8
$ws.Range("D32").Value = 290.35
$ws.Range("E32").Value = 1492.65
$ws.Range("F32").Value = 1091.95

$ws.Range("B33").Value = 130.3
$ws.Range("C33").Value = 89.55
$ws.Range("D33").Value = 300.7
$ws.Range("E33").Value = 1547.7
$ws.Range("F33").Value = 1130.7

$ws.Range("B34").Value = 134.35
$ws.Range("C34").Value = 91.65000000000001
$ws.Range("D34").Value = 309.45
$ws.Range("E34").Value = 1571.85
$ws.Range("F34").Value = 1167.45

$ws.Range("B35").Value = 134.95
$ws.Range("C35").Value = 93.05
$ws.Range("D35").Value = 315.05
$ws.Range("E35").Value = 1630.4
$ws.Range("F35").Value = 1198.3

$ws.Range("B36").Value = 139.25
$ws.Range("C36").Value = 96.05
$ws.Range("D36").Value = 326.2
$ws.Range("E36").Value = 1661.65
$ws.Range("F36").Value = 1243.25

$ws.Range("B37").Value = 138.4
$ws.Range("C37").Value = 97.8
$ws.Range("D37").Value = 332.9
$ws.Range("E37").Value = 1712.05
$ws.Range("F37").Value = 1275.9

$ws.Range("B38").Value = 143.65
$ws.Range("C38").Value = 103.1
$ws.Range("D38").Value = 347.2
$ws.Range("E38").Value = 1765.5
$ws.Range("F38").Value = 1315.8

$ws.Range("B39").Value = 146.6
$ws.Range("C39").Value = 110.35
$ws.Range("D39").Value = 350.35
$ws.Range("E39").Value = 1792.3
$ws.Range("F39").Value = 1348.9

$ws.Range("B40").Value = 147.85
$ws.Range("C40").Value = 117.25
$ws.Range("D40").Value = 358.2
$ws.Range("E40").Value = 1814.95
$ws.Range("F40").Value = 1380.35

$ws.Range("B41").Value = 153.2
$ws.Range("C41").Value = 123.4
$ws.Range("D41").Value = 370.1
$ws.Range("E41").Value = 1883.45
$ws.Range("F41").Value = 1425.7

